$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Split the run "()) Тогда //Инициализируем" so that an empty
#    string literal argument "" is inserted right after the opening
#    parenthesis: "()) ..." -> "("")) ...".  The three resulting
#    spans ("(", "\"\"" and ")) Тогда //Инициализируем") must stay as
#    three distinct runs (matching identical formatting), so a
#    transient formatting toggle is used on the inserted span to stop
#    the engine from re-coalescing it with its neighbours.
# ------------------------------------------------------------------
$rngText = $d.Content
$rngText.Find.Execute("()) Тогда //Инициализируем") | Out-Null
$splitPos = $rngText.Start + 1
$insertRng = $d.Range($splitPos, $splitPos)
$insertRng.InsertAfter('""')

$newSpan = $d.Range($splitPos, $splitPos + 2)
$newSpan.Font.Bold = 1
$newSpan.Font.Bold = 0

# ------------------------------------------------------------------
# 2) Move the "_GoBack" bookmark so that it spans from the start of
#    the "Попытка" paragraph through the end of the
#    " //Удаляем компоненту" paragraph (this also removes the old
#    bookmark that used to sit at the end of the document, because
#    Word only allows a single bookmark per name).
# ------------------------------------------------------------------
$rngStart = $d.Content
$rngStart.Find.Execute("Попытка") | Out-Null
$bmStart = $rngStart.Start

$rngEnd = $d.Content
$rngEnd.Find.Execute("//Удаляем компоненту") | Out-Null
$bmEnd = $rngEnd.Paragraphs.Item(1).Range.End

$bmRange = $d.Range($bmStart, $bmEnd)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
